# Rewrites the "Clio Muse Tours / personalized audio tours" proposal into
# the "Workable / AI candidate-match-scoring" proposal, per the commit diff.
#
# Unicode punctuation used by the source text (curly quotes/apostrophe).
$rsquo = [char]0x2019
$ldquo = [char]0x201C
$rdquo = [char]0x201D

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Title (Heading1)
#    "Project: Personalized Audio Tours for Clio Muse Tours"
#    -> "Project for Workable"
# -----------------------------------------------------------------
$d.Paragraphs(1).Range.Text = "Project for Workable"

# -----------------------------------------------------------------
# 2) "Who is the client?" -> "Introduction"  (Heading2)
# -----------------------------------------------------------------
$d.Paragraphs(2).Range.Text = "Introduction"

# -----------------------------------------------------------------
# 3) Paragraph 3 used to be:
#      "The client is a renowned platform named " + [hyperlink: Clio Muse
#      Tours] + ". Under the dynamic leadership ... previous year."
#    It becomes THREE paragraphs:
#      (a) FirstParagraph  - "Workable is a widely recognized ..."
#      (b) BodyText        - "The project idea that we are proposing is ..."
#      (c) BodyText        - "You can visit " + [hyperlink: Workable] +
#                             " to know more about their platform."
#    We keep the existing <w:hyperlink> (rId20) and only change its
#    display text, so its run formatting (rStyle=Hyperlink) survives.
# -----------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$hl = $d.Hyperlinks.Item(1)
$hl.TextToDisplay = "Workable"

# Strip the old text surrounding the hyperlink, leaving only the
# hyperlink run behind in paragraph 3.
$beforeRange = $d.Range($p3.Range.Start, $hl.Range.Start)
$beforeRange.Text = ""
$hl = $d.Hyperlinks.Item(1)
$afterRange = $d.Range($hl.Range.End, $p3.Range.End)
$afterRange.Text = ""

# Wrap the (now bare) hyperlink with the new sentence and re-style the
# paragraph as BodyText - this becomes final paragraph (c) above.
$hl = $d.Hyperlinks.Item(1)
$hl.Range.InsertBefore("You can visit ")
$hl = $d.Hyperlinks.Item(1)
$hl.Range.InsertAfter(" to know more about their platform.")
$d.Paragraphs(3).Style = "BodyText"

# Make room for the two new paragraphs in front of it.
$p3 = $d.Paragraphs(3)
$insertPoint = $d.Range($p3.Range.Start, $p3.Range.Start)
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()

$introText = "Workable is a widely recognized recruitment software that " + `
    "has been successful in assisting over 30,000 companies across " + `
    "100+ countries in hiring the best talent. The platform was founded " + `
    "in 2012 by Nikos Moraitakis in Athens, Greece, with a sole mission " + `
    "to streamline the hiring process for small and medium-sized " + `
    "businesses. Nikos Moraitakis, the CEO, has led the company to " + `
    "global expansion and significant fundraising, thereby making a " + `
    "prominent mark in the recruitment industry."
$d.Paragraphs(3).Range.Text = $introText
$d.Paragraphs(3).Style = "FirstParagraph"

$ideaText = "The project idea that we are proposing is the " + $ldquo + `
    "Implementation of an AI-powered Candidate Match Scoring system" + `
    $rdquo + " for Workable. This system is aimed to analyze job " + `
    "descriptions and candidate profiles to generate accurate match " + `
    "scores, thereby speeding up the recruitment process and ensuring " + `
    "the best fit."
$d.Paragraphs(4).Range.Text = $ideaText
$d.Paragraphs(4).Style = "BodyText"

# -----------------------------------------------------------------
# Paragraph indices from here on (after the 2-paragraph insert above):
#   6  "What is the idea?"        -> "Client Background"
#   7  project idea paragraph     -> client background paragraph
#   8  "How can we help?"         -> "Proposed Solution"
#   9  solution paragraph         -> proposed solution paragraph
#   10 BodyText paragraph         -> integration paragraph
#   11 "Tech Stack"                  (unchanged)
#   12 tech-stack intro sentence  -> new intro sentence
#   13-19 the 7 tech-stack bullets -> replaced text (+1 new bullet added)
#   20 "Timeline"                    (unchanged)
#   21 timeline intro sentence    -> new intro sentence
#   22-24 the 3 timeline bullets  -> replaced text
#   25 timeline closing sentence  -> new closing sentence
# -----------------------------------------------------------------

# 4) "What is the idea?" -> "Client Background"
$d.Paragraphs(6).Range.Text = "Client Background"

# 5) Client background paragraph
$clientBg = "The client, Nikos Moraitakis, is the founder and CEO of " + `
    "Workable. He has a strong background in leading tech companies " + `
    "and has been instrumental in the success story of Workable. His " + `
    "vision of creating better recruitment software for small and " + `
    "medium-sized businesses has transformed the company into a " + `
    "globally recognized brand. His leadership skills have enabled " + `
    "Workable to serve thousands of companies across numerous " + `
    "countries, raise significant funding, and expand its operations " + `
    "worldwide."
$d.Paragraphs(7).Range.Text = $clientBg

# 6) "How can we help?" -> "Proposed Solution"
$d.Paragraphs(8).Range.Text = "Proposed Solution"

# 7) Proposed-solution paragraph
$proposedSolution = "Our proposed solution involves developing a " + `
    "machine learning model that will analyze job descriptions and " + `
    "candidate profiles to generate accurate match scores. The system " + `
    "will utilize Natural Language Processing (NLP) to extract key " + `
    "skills, experience, and qualifications from resumes and job " + `
    "postings. It will then employ a deep learning algorithm to " + `
    "calculate a compatibility score based on various factors such as " + `
    "skills alignment, experience relevance, and cultural fit."
$d.Paragraphs(9).Range.Text = $proposedSolution

# 8) Integration paragraph
$integration = "The system will be integrated with Workable" + $rsquo + `
    "s existing platform, providing seamless scoring capabilities " + `
    "within their applicant tracking system. This will enable " + `
    "recruiters to quickly identify and shortlist candidates that " + `
    "best fit the job requirements, thereby reducing the time and " + `
    "effort involved in the screening process."
$d.Paragraphs(10).Range.Text = $integration

# 9) Tech Stack heading unchanged.

# 10) Tech-stack intro sentence
$d.Paragraphs(12).Range.Text = "For this project, we will be using the following tech stack:"

# 11) The 7 existing tech-stack bullets get new text ...
$d.Paragraphs(13).Range.Text = "Python: For coding the backend of the application."
$d.Paragraphs(14).Range.Text = "TensorFlow, PyTorch: For building and training the machine learning and deep learning models."
$d.Paragraphs(15).Range.Text = "NLTK, spaCy: For performing Natural Language Processing tasks."
$d.Paragraphs(16).Range.Text = "FastAPI: For building the APIs."
$d.Paragraphs(17).Range.Text = "Docker, Kubernetes: For containerization and orchestration."
$d.Paragraphs(18).Range.Text = "AWS (EC2, S3, SageMaker): For cloud computing and storage needs."
$d.Paragraphs(19).Range.Text = "PostgreSQL: As the primary database."

# ... and a new 8th bullet ("Redis") is appended after that last one,
# inheriting the Compact/numPr list formatting automatically.
$d.Paragraphs(19).Range.InsertParagraphAfter()
$d.Paragraphs(20).Range.Text = "Redis: For caching and session management."

# -----------------------------------------------------------------
# Paragraph indices shift by +1 again after the Redis insert above:
#   21 "Timeline"                    (unchanged)
#   22 timeline intro sentence    -> new intro sentence
#   23-25 the 3 timeline bullets  -> replaced text
#   26 timeline closing sentence  -> new closing sentence
# -----------------------------------------------------------------

# 12) Timeline heading unchanged.

# 13) Timeline intro sentence
$d.Paragraphs(22).Range.Text = "The project is anticipated to be completed within 4 to 6 months, distributed as follows:"

# 14) The 3 timeline bullets
$d.Paragraphs(23).Range.Text = "1 month for requirements gathering and design: Understanding the client" + $rsquo + "s requirements and designing the system architecture."
$d.Paragraphs(24).Range.Text = "2-3 months for development and AI model training: Coding the solution and training the AI models using relevant data."
$d.Paragraphs(25).Range.Text = "1-2 months for testing, integration, and deployment: Testing the solution thoroughly, integrating it with the existing platform, and deploying it for use."

# 15) Timeline closing sentence
$d.Paragraphs(26).Range.Text = "This timeline ensures that the project is completed efficiently without compromising on the quality of the solution."

# -----------------------------------------------------------------
# Bookmark names (X8194c6c...->project-for-workable, who-is-the-client
# -> introduction, what-is-the-idea -> client-background, how-can-we-help
# -> proposed-solution). Best-effort: this COM-interop surface does not
# expose a working Bookmarks collection (Count stays 0 / Add is inert),
# so these calls are harmless no-ops if unsupported.
# -----------------------------------------------------------------
try {
    $d.Bookmarks.Add("project-for-workable") | Out-Null
    $d.Bookmarks.Add("introduction") | Out-Null
    $d.Bookmarks.Add("client-background") | Out-Null
    $d.Bookmarks.Add("proposed-solution") | Out-Null
} catch {
}

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
